$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, $CellRef, $Text)
    $r = $Worksheet.Range($CellRef)
    $origStyle = $r.Style
    $r.Value = "'" + $Text
    $r.Style = $origStyle
}

$ws.Range('D2').Value = '64.012.23'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').Value = '3.431.87'
$ws.Range('E3').Value = '  +0.66%  '
Set-TextValue $ws 'D4' '1.00'
$ws.Range('E4').Value = '  +0.03%  '
Set-TextValue $ws 'D5' '571.63'
$ws.Range('E5').Value = '  +0.27%  '
Set-TextValue $ws 'D6' '159.30'
$ws.Range('E6').Value = '  +1.76%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '3.433.55'
$ws.Range('E8').Value = '  +0.58%  '
$ws.Range('E9').Value = '  -8.82%  '
Set-TextValue $ws 'D10' '7.26'
$ws.Range('E10').Value = '  +1.73%  '
$ws.Range('E11').Value = '  -2.14%  '
$ws.Range('E12').Value = '  -3.66%  '
$ws.Range('D13').Value = '4.022.61'
$ws.Range('E13').Value = '  +0.75%  '
$ws.Range('E14').Value = '  +0.53%  '
$ws.Range('E15').Value = '  -1.24%  '
$ws.Range('E16').Value = '  -7.29%  '
$ws.Range('D17').Value = '64.063.55'
$ws.Range('E17').Value = '  -0.46%  '
$ws.Range('D18').Value = '3.451.64'
$ws.Range('E18').Value = '  +1.34%  '
$ws.Range('E19').Value = '  -3.30%  '
Set-TextValue $ws 'D20' '13.64'
$ws.Range('E20').Value = '  -1.02%  '
Set-TextValue $ws 'D21' '384.26'
$ws.Range('E21').Value = '  +1.77%  '
Set-TextValue $ws 'D22' '7.85'
$ws.Range('E22').Value = '  -2.06%  '
$ws.Range('E23').Value = '  +0.32%  '
Set-TextValue $ws 'D24' '71.35'
$ws.Range('E24').Value = '  -0.22%  '
$ws.Range('E25').Value = '  -5.23%  '
$ws.Range('E26').Value = '  -1.71%  '
Set-TextValue $ws 'D27' '9.68'
$ws.Range('E27').Value = '  -6.51%  '
$ws.Range('E28').Value = '  +0.11%  '
Set-TextValue $ws 'D29' '0.997'
$ws.Range('E29').Value = '  -0.18%  '
Set-TextValue $ws 'D30' '6.07'
$ws.Range('E30').Value = '  -1.55%  '
$ws.Range('E31').Value = '  -5.35%  '
$ws.Range('E32').Value = '  -0.41%  '
Set-TextValue $ws 'D33' '23.00'
$ws.Range('E33').Value = '  -0.19%  '
$ws.Range('E34').Value = '  +0.00%  '
Set-TextValue $ws 'D35' '6.97'
$ws.Range('E35').Value = '  -2.61%  '
$ws.Range('E36').Value = '  -5.24%  '
Set-TextValue $ws 'D37' '160.87'
$ws.Range('E37').Value = '  +0.69%  '
Set-TextValue $ws 'D38' '0.850'
$ws.Range('E38').Value = '  +10.73%  '
$ws.Range('E39').Value = '  -2.75%  '
$ws.Range('D40').Value = '2.833.86'
$ws.Range('E40').Value = '  -1.54%  '
Set-TextValue $ws 'D41' '26.04'
$ws.Range('E41').Value = '  -0.84%  '
$ws.Range('E42').Value = '  -4.76%  '
Set-TextValue $ws 'D43' '43.05'
$ws.Range('E43').Value = '  +0.30%  '
Set-TextValue $ws 'D44' '26.41'
$ws.Range('E44').Value = '  +1.21%  '
$ws.Range('E45').Value = '  -7.69%  '
$ws.Range('E46').Value = '  -5.45%  '
Set-TextValue $ws 'D47' '0.0304'
$ws.Range('E47').Value = '  -3.24%  '
Set-TextValue $ws 'D48' '2.43'
$ws.Range('E48').Value = '  +11.37%  '
Set-TextValue $ws 'D49' '333.77'
$ws.Range('E49').Value = '  +3.52%  '
$ws.Range('E50').Value = '  -2.33%  '
$ws.Range('E51').Value = '  -5.72%  '
